# Update FlashScore odds data for 2024-11-16 workbook (Jogos_da_Semana)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: updated odds values
$ws.Range("G2").Value = 1.8
$ws.Range("I2").Value = 4.5
$ws.Range("J2").Value = 2.4
$ws.Range("L2").Value = 4.5
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 1.95
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 2.75
$ws.Range("W2").Value = 7.5
$ws.Range("AB2").Value = 26
$ws.Range("AE2").Value = 15
$ws.Range("AF2").Value = 51
$ws.Range("AI2").Value = 23
$ws.Range("AN2").Value = 3.75
$ws.Range("AO2").Value = 9.5
$ws.Range("AS2").Value = 151
$ws.Range("AT2").Value = 2.75
$ws.Range("AX2").Value = 23
$ws.Range("AZ2").Value = 81
$ws.Range("BA2").Value = 101

# Row 3: updated odds values
$ws.Range("G3").Value = 4.2
$ws.Range("I3").Value = 1.95
$ws.Range("J3").Value = 4.5
$ws.Range("Y3").Value = 15
$ws.Range("AE3").Value = 17
$ws.Range("AI3").Value = 8.5
$ws.Range("AO3").Value = 23
$ws.Range("AQ3").Value = 81
$ws.Range("AX3").Value = 11

# Row 5: updated odds values
$ws.Range("G5").Value = 2.12
$ws.Range("H5").Value = 3.2
$ws.Range("I5").Value = 3.15
$ws.Range("J5").Value = 2.77
$ws.Range("K5").Value = 2.12
$ws.Range("L5").Value = 3.65
$ws.Range("N5").Value = 7.7
$ws.Range("S5").Value = 1.37
$ws.Range("T5").Value = 2.85
$ws.Range("X5").Value = 11.25
$ws.Range("Y5").Value = 8.5
$ws.Range("Z5").Value = 21
$ws.Range("AA5").Value = 16.5
$ws.Range("AB5").Value = 23
$ws.Range("AC5").Value = 7.7
$ws.Range("AD5").Value = 6.4
$ws.Range("AI5").Value = 19.5
$ws.Range("AJ5").Value = 10.75
$ws.Range("AK5").Value = 45
$ws.Range("AL5").Value = 25
$ws.Range("AM5").Value = 28
$ws.Range("AN5").Value = 4.2
$ws.Range("AO5").Value = 11.5
$ws.Range("AP5").Value = 18
$ws.Range("AQ5").Value = 45
$ws.Range("AR5").Value = 70
$ws.Range("AT5").Value = 2.85
$ws.Range("AW5").Value = 5.3
$ws.Range("AX5").Value = 17
$ws.Range("AY5").Value = 21
$ws.Range("AZ5").Value = 80
$ws.Range("BA5").Value = 100

# Remove the Venezuela LIGA FUTVE match (was row 7); remaining rows shift up
# and the sheet dimension is recalculated automatically (A1:BD7 -> A1:BD6)
$ws.Rows.Item(7).Delete()

